
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column at E for German (shifts old E,F -> F,G) ---
$ws.Columns.Item(5).Insert()

# --- 2. Swap columns A and B (A held images, B held EN text; now reversed) ---
$ws.Range("Z2:Z4").Clear()
$ws.Range("A2:A4").Copy($ws.Range("Z2:Z4"))
$ws.Range("B2:B4").Copy($ws.Range("A2:A4"))
$ws.Range("Z2:Z4").Copy($ws.Range("B2:B4"))
$ws.Range("Z2:Z4").Clear()

# --- 3. Relabel header row 1 ---
$ws.Range("A1").Value = 'inst_pics'
$ws.Range("B1").Value = 'inst_msg_EN'
$ws.Range("C1").Value = 'inst_msg_ES'
$ws.Range("D1").Value = 'inst_msg_FR'
$ws.Range("E1").Value = 'inst_msg_DE'
$ws.Range("F1").Value = 'image_w'
$ws.Range("G1").Value = 'image_h'

# --- 4. Fill in new German column content ---
$ws.Range("E2").Value = "Anweisungen:`nIn dieser Aufgabe sehen Sie Bilder von linken oder rechten Händen, deren Handflächen nach oben oder unten zeigen. Die Bilder werden in verschiedenen Winkeln gedreht.`nIhre Aufgabe ist es zu bestimmen, ob das Bild einer linken oder rechten Hand entspricht.`nIhr Ziel ist es, sowohl SCHNELL als auch GENAU zu antworten.`nJedes Bild wird angezeigt, bis Sie geantwortet haben. Das nächste Bild erscheint automatisch."
$ws.Range("E3").Value = "Bitte verwenden Sie nur Ihren ZEIGEFINGER und MITTELFINGER der RECHTEN HAND, um zu antworten.`nLegen Sie Ihren Zeigefinger auf die Taste „G“ und den Mittelfinger auf die Taste „H“ Ihrer Tastatur.`nZum Antworten:`nLinke Hand = G | H = Rechte Hand`nSie müssen Ihre Hand während der gesamten Aufgabe auf der Tastatur halten.`nHalten Sie Ihre andere Hand auf dem Tisch, in derselben Position und so ruhig wie möglich."
$ws.Range("E4").Value = "Nach jedem Bild erhalten Sie ein kurzes Feedback zu Ihrer Antwort:`nWenn Sie korrekt antworten, wird das entsprechende Feld grün.`nWenn Sie falsch antworten, wird das entsprechende Feld rot.`nDenken Sie daran, dass Ihr Ziel darin besteht, so genau und schnell wie möglich zu antworten."

# --- 5. Column widths ---
$ws.Columns.Item(1).ColumnWidth = 33.90625
$ws.Columns.Item(2).ColumnWidth = 16.81640625
$ws.Columns.Item(3).ColumnWidth = 14.7265625
$ws.Columns.Item(4).ColumnWidth = 15.36328125
$ws.Columns.Item(5).ColumnWidth = 15.36328125
$ws.Columns.Item(6).ColumnWidth = 8.7265625
$ws.Columns.Item(7).ColumnWidth = 8.36328125

# --- 6. Header row style: blue fill + white font (set cleanly on A1, then propagate via format-only paste) ---
$ws.Range("A1").Font.Color = 16777215
$ws.Range("A1").Interior.Color = 15773696
$ws.Range("A1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 7. Row heights back to 24.5 for data rows ---
$ws.Rows.Item(2).RowHeight = 24.5
$ws.Rows.Item(3).RowHeight = 24.5
$ws.Rows.Item(4).RowHeight = 24.5

# --- 8. Dimension / Selection ---
$ws.Range("E4").Select()
